# Auto-generated edit script for VerveStacks ZAF workbook update
# Updates 'efficiency'/'ncap_cost'/'ncap_fom'/'act_cost'/'AF' columns on ccs_retrofits sheet
# and 'ncap_fom' column on existing_stock sheet, per commit 'Updated ZAF model - 2025-07-29 12:13'

$wb = $excel.ActiveWorkbook

$wsCcs = $wb.Worksheets.Item("ccs_retrofits")
$wsExisting = $wb.Worksheets.Item("existing_stock")

# --- ccs_retrofits sheet (sheet2): columns D (efficiency), E (ncap_cost), F (ncap_fom), G (act_cost), H (AF) ---
$wsCcs.Range("D4:D7").Value = 0.07441920000000003
$wsCcs.Range("D8:D9").Value = 0.16989120000000002
$wsCcs.Range("D10:D17").Value = 0.06566400000000001
$wsCcs.Range("D18:D19").Value = 0.16989120000000002
$wsCcs.Range("D20:D23").Value = 0.17488800000000002
$wsCcs.Range("D24").Value = 0.07004160000000002
$wsCcs.Range("D25:D27").Value = 0.06566400000000001
$wsCcs.Range("D28:D32").Value = 0.07004160000000002
$wsCcs.Range("D33:D34").Value = 0.06566400000000001
$wsCcs.Range("D35:D39").Value = 0.07004160000000002
$wsCcs.Range("D40:D46").Value = 0.06566400000000001
$wsCcs.Range("D47:D49").Value = 0.17488800000000002
$wsCcs.Range("D50:D52").Value = 0.17988480000000004
$wsCcs.Range("D53:D58").Value = 0.16989120000000002
$wsCcs.Range("D59:D64").Value = 0.17488800000000002
$wsCcs.Range("D65:D67").Value = 0.17988480000000004
$wsCcs.Range("D68:D70").Value = 0.20839680000000005
$wsCcs.Range("D71:D75").Value = 0.17488800000000002
$wsCcs.Range("D76").Value = 0.17988480000000004
$wsCcs.Range("D77:D79").Value = 0.16989120000000002
$wsCcs.Range("D80:D82").Value = 0.17488800000000002
$wsCcs.Range("D83:D87").Value = 0.06566400000000001
$wsCcs.Range("D88:D93").Value = 0.17488800000000002
$wsCcs.Range("D94").Value = 0.059097600000000014
$wsCcs.Range("D95").Value = 0.07004160000000004
$wsCcs.Range("D96:D97").Value = 0.2402352
$wsCcs.Range("D98:D102").Value = 0.24602400000000008
$wsCcs.Range("D103:D107").Value = 0.2402352
$wsCcs.Range("E20:E23").Value = 2445.0
$wsCcs.Range("E47:E52").Value = 2445.0
$wsCcs.Range("E59:E67").Value = 2445.0
$wsCcs.Range("E68:E70").Value = 2126.0
$wsCcs.Range("E71:E76").Value = 2445.0
$wsCcs.Range("E80:E82").Value = 2445.0
$wsCcs.Range("E88:E93").Value = 2445.0
$wsCcs.Range("E96:E107").Value = 2126.0
$wsCcs.Range("F20:F23").Value = 30.4
$wsCcs.Range("F47:F52").Value = 30.4
$wsCcs.Range("F59:F67").Value = 30.4
$wsCcs.Range("F68:F70").Value = 25.1
$wsCcs.Range("F71:F76").Value = 30.4
$wsCcs.Range("F80:F82").Value = 30.4
$wsCcs.Range("F88:F93").Value = 30.4
$wsCcs.Range("F96:F107").Value = 25.1
$wsCcs.Range("G20:G23").Value = 3.52
$wsCcs.Range("G47:G52").Value = 3.52
$wsCcs.Range("G59:G67").Value = 3.52
$wsCcs.Range("G68:G70").Value = 3.22
$wsCcs.Range("G71:G76").Value = 3.52
$wsCcs.Range("G80:G82").Value = 3.52
$wsCcs.Range("G88:G93").Value = 3.52
$wsCcs.Range("G96:G107").Value = 3.22
$wsCcs.Range("H20:H23").Value = 0.7277
$wsCcs.Range("H47:H52").Value = 0.7277
$wsCcs.Range("H59:H67").Value = 0.7277
$wsCcs.Range("H68:H70").Value = 0.7941999999999999
$wsCcs.Range("H71:H76").Value = 0.7277
$wsCcs.Range("H80:H82").Value = 0.7277
$wsCcs.Range("H88:H93").Value = 0.7277
$wsCcs.Range("H96:H107").Value = 0.7941999999999999

# --- existing_stock sheet (sheet3): column F (ncap_fom) ---
$wsExisting.Range("F4").Value = 0.23040000000000008
$wsExisting.Range("F5:F6").Value = 0.24480000000000005
$wsExisting.Range("F7:F8").Value = 0.24480000000000007
$wsExisting.Range("F9:F10").Value = 0.24480000000000005
$wsExisting.Range("F11:F18").Value = 0.21600000000000003
$wsExisting.Range("F19:F20").Value = 0.24480000000000005
$wsExisting.Range("F21:F24").Value = 0.252
$wsExisting.Range("F25").Value = 0.23040000000000005
$wsExisting.Range("F26:F28").Value = 0.21600000000000003
$wsExisting.Range("F29:F33").Value = 0.23040000000000005
$wsExisting.Range("F34:F35").Value = 0.21600000000000003
$wsExisting.Range("F36:F40").Value = 0.23040000000000005
$wsExisting.Range("F41:F47").Value = 0.21600000000000003
$wsExisting.Range("F48:F50").Value = 0.252
$wsExisting.Range("F51:F53").Value = 0.25920000000000004
$wsExisting.Range("F54:F59").Value = 0.24480000000000005
$wsExisting.Range("F60:F65").Value = 0.252
$wsExisting.Range("F66:F71").Value = 0.25920000000000004
$wsExisting.Range("F72:F76").Value = 0.252
$wsExisting.Range("F77").Value = 0.25920000000000004
$wsExisting.Range("F78:F80").Value = 0.24480000000000005
$wsExisting.Range("F81:F83").Value = 0.252
$wsExisting.Range("F84:F88").Value = 0.21600000000000003
$wsExisting.Range("F89:F94").Value = 0.252
$wsExisting.Range("F95").Value = 0.19440000000000002
$wsExisting.Range("F96:F97").Value = 0.2988
$wsExisting.Range("F98:F102").Value = 0.30600000000000005
$wsExisting.Range("F103:F107").Value = 0.2988
